# CreateOrder-Event.xlsx -- release 1.6.3 update
#
# The single content/authoring change in this revision is on the
# "OrderCreated-Event" sheet: the MessageType value in row 3 (cell K3)
# is renamed from "JSONMessageType" to "JSONType". Saving also records
# the author's last-active cell/selection as K4 (one row below the
# edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrderCreated-Event")

# K1 = "MessageType" header; K3 holds the row-3 message type value.
$ws.Range("K3").Value = "JSONType"

# Record the cursor/selection left on the sheet after the edit.
$ws.Range("K4").Select()
